# Update "Förändrad" date (column C) for rows 2-9 from 2023-11-03 to 2023-11-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = (Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0)
}
